$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 16:30"

# Row 4
$ws.Range("B4").Value = 7006994
$ws.Range("C4").Value = 2226
$ws.Range("D4").Value = 4250545
$ws.Range("E4").Value = 2552311
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 204138

# Row 5
$ws.Range("B5").Value = 5517596
$ws.Range("C5").Value = 31984
$ws.Range("D5").Value = 4428913
$ws.Range("E5").Value = 1000452
$ws.Range("G5").Value = 322
$ws.Range("H5").Value = 88231

# Row 20
$ws.Range("B20").Value = 322856
$ws.Range("C20").Value = 3821
$ws.Range("D20").Value = 258075
$ws.Range("E20").Value = 56156
$ws.Range("G20").Value = 70
$ws.Range("H20").Value = 8625

# Row 51
$ws.Range("A51").Value = "Portugal"
$ws.Range("B51").Value = 69200
$ws.Range("C51").Value = 623
$ws.Range("D51").Value = 45736
$ws.Range("E51").Value = 21544
$ws.Range("G51").Value = 8
$ws.Range("H51").Value = 1920

# Row 52
$ws.Range("A52").Value = "Etiopia"
$ws.Range("B52").Value = 68820
$ws.Range("D52").Value = 28314
$ws.Range("E52").Value = 39410
$ws.Range("H52").Value = 1096

# Row 64
$ws.Range("B64").Value = 46796
$ws.Range("C64").Value = 200
$ws.Range("D64").Value = 35018
$ws.Range("E64").Value = 10567
$ws.Range("G64").Value = 8
$ws.Range("H64").Value = 1211

# Row 73
$ws.Range("A73").Value = "Serbia"
$ws.Range("B73").Value = 32938
$ws.Range("C73").Value = 30
$ws.Range("D73").Value = 31536
$ws.Range("E73").Value = 659
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 743

# Row 74
$ws.Range("A74").Value = "Irlanda"
$ws.Range("B74").Value = 32933
$ws.Range("D74").Value = 23364
$ws.Range("E74").Value = 7777
$ws.Range("H74").Value = 1792

# Row 86
$ws.Range("B86").Value = 16780
$ws.Range("C86").Value = 45
$ws.Range("D86").Value = 13949
$ws.Range("E86").Value = 2131
$ws.Range("G86").Value = 7
$ws.Range("H86").Value = 700

# Row 93
$ws.Range("B93").Value = 12954
$ws.Range("C93").Value = 57
$ws.Range("E93").Value = 2316

# Row 96
$ws.Range("A96").Value = "Namibia"
$ws.Range("B96").Value = 10526
$ws.Range("C96").Value = 149
$ws.Range("D96").Value = 8112
$ws.Range("E96").Value = 2301
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 113

# Row 97
$ws.Range("A97").Value = "Consejo Danes para los Refugiados"
$ws.Range("B97").Value = 10519
$ws.Range("C97").Value = 4
$ws.Range("D97").Value = 9952
$ws.Range("E97").Value = 296
$ws.Range("H97").Value = 271

# Row 105
$ws.Range("B105").Value = 8624
$ws.Range("C105").Value = 5
$ws.Range("D105").Value = 6482
$ws.Range("E105").Value = 1921

# Row 118
$ws.Range("B118").Value = 5141
$ws.Range("C118").Value = 50
$ws.Range("D118").Value = 4462
$ws.Range("E118").Value = 563
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 116

# Row 138
$ws.Range("B138").Value = 3465
$ws.Range("C138").Value = 23
$ws.Range("E138").Value = 490

# Row 140
$ws.Range("B140").Value = 3370
$ws.Range("C140").Value = 55
$ws.Range("E140").Value = 1607

# Row 141
$ws.Range("B141").Value = 3290
$ws.Range("C141").Value = 3
$ws.Range("E141").Value = 177

# Row 148
$ws.Range("B148").Value = 2377
$ws.Range("C148").Value = 31
$ws.Range("D148").Value = 2125
$ws.Range("E148").Value = 242

# Row 159
$ws.Range("D159").Value = 1369
$ws.Range("E159").Value = 209

# Row 163
$ws.Range("B163").Value = 1336
$ws.Range("C163").Value = 1
$ws.Range("D163").Value = 1218
$ws.Range("E163").Value = 36

# Row 204
$ws.Range("A204").Value = "Santa Lucia"

# Row 205
$ws.Range("A205").Value = "Timor Oriental"

# Row 214
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Row 215
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

Write-Output "Applied paises.xlsx update"